$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 263, shifting existing rows 263-293 down to 266-296
$ws.Range("A263:T265").Insert()

# Populate the 3 newly inserted rows with the new week of Chirimoya price data (2023-07-24)
$ws.Range("A263").Value = 3
$ws.Range("B263").Value = 'Femacal de La Calera'
$ws.Range("C263").Value = 'Coquimbo'
$ws.Range("D263").Value = 45131
$ws.Range("E263").Value = 5
$ws.Range("F263").Value = 'Fruta'
$ws.Range("G263").Value = 100107
$ws.Range("H263").Value = 'Otros'
$ws.Range("I263").Value = 100107002
$ws.Range("J263").Value = 'Chirimoya'
$ws.Range("K263").Value = 'Cultivar IV Región'
$ws.Range("L263").Value = 'Especial'
$ws.Range("M263").Value = 56
$ws.Range("N263").Value = 30000
$ws.Range("O263").Value = 30000
$ws.Range("P263").Value = 30000
$ws.Range("Q263").Value = '$/bandeja 10 kilos'
$ws.Range("R263").Value = 'Provincia del Elquí'
$ws.Range("S263").Value = 3000
$ws.Range("T263").Value = 10

$ws.Range("A264").Value = 3
$ws.Range("B264").Value = 'Femacal de La Calera'
$ws.Range("C264").Value = 'Coquimbo'
$ws.Range("D264").Value = 45131
$ws.Range("E264").Value = 5
$ws.Range("F264").Value = 'Fruta'
$ws.Range("G264").Value = 100107
$ws.Range("H264").Value = 'Otros'
$ws.Range("I264").Value = 100107002
$ws.Range("J264").Value = 'Chirimoya'
$ws.Range("K264").Value = 'Cultivar IV Región'
$ws.Range("L264").Value = 'Primera'
$ws.Range("M264").Value = 60
$ws.Range("N264").Value = 28000
$ws.Range("O264").Value = 28000
$ws.Range("P264").Value = 28000
$ws.Range("Q264").Value = '$/bandeja 10 kilos'
$ws.Range("R264").Value = 'Provincia del Elquí'
$ws.Range("S264").Value = 2800
$ws.Range("T264").Value = 10

$ws.Range("A265").Value = 3
$ws.Range("B265").Value = 'Femacal de La Calera'
$ws.Range("C265").Value = 'Coquimbo'
$ws.Range("D265").Value = 45131
$ws.Range("E265").Value = 5
$ws.Range("F265").Value = 'Fruta'
$ws.Range("G265").Value = 100107
$ws.Range("H265").Value = 'Otros'
$ws.Range("I265").Value = 100107002
$ws.Range("J265").Value = 'Chirimoya'
$ws.Range("K265").Value = 'Cultivar IV Región'
$ws.Range("L265").Value = 'Segunda'
$ws.Range("M265").Value = 36
$ws.Range("N265").Value = 25000
$ws.Range("O265").Value = 25000
$ws.Range("P265").Value = 25000
$ws.Range("Q265").Value = '$/bandeja 10 kilos'
$ws.Range("R265").Value = 'Provincia del Elquí'
$ws.Range("S265").Value = 2500
$ws.Range("T265").Value = 10
